# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly-crawled counts, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" - rows 2-7
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14436
$wsExhibit.Range("F3").Value = 339
$wsExhibit.Range("F4").Value = 694
$wsExhibit.Range("F5").Value = 241
$wsExhibit.Range("F6").Value = 577
$wsExhibit.Range("F7").Value = 1516

# Sheet "全部类型" - rows 2-5, 8, 9
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14436
$wsAll.Range("F3").Value = 339
$wsAll.Range("F4").Value = 694
$wsAll.Range("F5").Value = 241
$wsAll.Range("F8").Value = 577
$wsAll.Range("F9").Value = 1516
